$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: "Model Accuracy (-0.6, 0.6, 0.6)"
# Add new columns C:G (Market threshold, Market min, Market max, Recall,
# Precision) with header formatting copied from the existing B1 header,
# and update/add the numeric data for rows 2-6.
# ----------------------------------------------------------------------
$wsAcc = $wb.Worksheets.Item("Model Accuracy (-0.6, 0.6, 0.6)")

# Copy the header style (bold, border, centered) from B1 onto the new
# header cells C1:G1 so they reuse the same cell style as B1.
$wsAcc.Range("B1").Copy()
$wsAcc.Range("C1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsAcc.Range("C1").Value = "Market threshold"
$wsAcc.Range("D1").Value = "Market min"
$wsAcc.Range("E1").Value = "Market max"
$wsAcc.Range("F1").Value = "Recall"
$wsAcc.Range("G1").Value = "Precision"

# Row 2 - TOTALENERGIES SE
$wsAcc.Range("B2").Value = 62.59168704156479
$wsAcc.Range("C2").Value = 0.05450546436368681
$wsAcc.Range("D2").Value = -15.55441
$wsAcc.Range("E2").Value = 15.06418
$wsAcc.Range("F2").Value = 0
$wsAcc.Range("G2").Value = 0

# Row 3 - FMC CORP
$wsAcc.Range("B3").Value = 37.16381418092909
$wsAcc.Range("C3").Value = 0.009583939973006913
$wsAcc.Range("D3").Value = -19.35264
$wsAcc.Range("E3").Value = 13.70093
$wsAcc.Range("F3").Value = 2.680965147453083
$wsAcc.Range("G3").Value = 25.64102564102564

# Row 4 - BP PLC
$wsAcc.Range("B4").Value = 92.1760391198044
$wsAcc.Range("C4").Value = 0.04158117063764853
$wsAcc.Range("D4").Value = -18.75314
$wsAcc.Range("E4").Value = 23.33066
$wsAcc.Range("F4").Value = 0
$wsAcc.Range("G4").Value = 0

# Row 5 - STORA ENSO
$wsAcc.Range("B5").Value = 81.11246943765281
$wsAcc.Range("C5").Value = 0.02983403801513819
$wsAcc.Range("D5").Value = -12.78028
$wsAcc.Range("E5").Value = 12.42348
$wsAcc.Range("F5").Value = 0
$wsAcc.Range("G5").Value = 0

# Row 6 - BHP GROUP
$wsAcc.Range("B6").Value = 94.8044009779951
$wsAcc.Range("C6").Value = 0.08368817696170747
$wsAcc.Range("D6").Value = -16.47904
$wsAcc.Range("E6").Value = 14.94325
$wsAcc.Range("F6").Value = 0
$wsAcc.Range("G6").Value = 0

# ----------------------------------------------------------------------
# Sheet 2: "Confusion Matrix TOTALENERGIES SE (-0.6, 0.6, 0.6)"
# ----------------------------------------------------------------------
$wsCm1 = $wb.Worksheets.Item("Confusion Matrix TOTALENERGIES SE (-0.6, 0.6, 0.6)")
$wsCm1.Range("B3").Value = 9
$wsCm1.Range("C3").Value = 1022
$wsCm1.Range("D3").Value = 8

# ----------------------------------------------------------------------
# Sheet 3: "Confusion Matrix FMC CORP (-0.6, 0.6, 0.6)"
# ----------------------------------------------------------------------
$wsCm2 = $wb.Worksheets.Item("Confusion Matrix FMC CORP (-0.6, 0.6, 0.6)")
$wsCm2.Range("B2").Value = 10
$wsCm2.Range("C2").Value = 22
$wsCm2.Range("D2").Value = 7

$wsCm2.Range("B3").Value = 315
$wsCm2.Range("C3").Value = 558
$wsCm2.Range("D3").Value = 309

$wsCm2.Range("B4").Value = 48
$wsCm2.Range("C4").Value = 72
$wsCm2.Range("D4").Value = 40

# ----------------------------------------------------------------------
# Sheet 4: "Confusion Matrix BP PLC (-0.6, 0.6, 0.6)"
# ----------------------------------------------------------------------
$wsCm3 = $wb.Worksheets.Item("Confusion Matrix BP PLC (-0.6, 0.6, 0.6)")
$wsCm3.Range("B3").Value = 39
$wsCm3.Range("C3").Value = 1507
$wsCm3.Range("D3").Value = 41

$wsCm3.Range("B4").Value = 1
$wsCm3.Range("C4").Value = 10

# ----------------------------------------------------------------------
# Sheet 5: "Confusion Matrix STORA ENSO (-0.6, 0.6, 0.6)"
# ----------------------------------------------------------------------
$wsCm4 = $wb.Worksheets.Item("Confusion Matrix STORA ENSO (-0.6, 0.6, 0.6)")
$wsCm4.Range("B3").Value = 108
$wsCm4.Range("C3").Value = 1327
$wsCm4.Range("D3").Value = 107

$wsCm4.Range("B4").Value = 2
$wsCm4.Range("C4").Value = 23

# ----------------------------------------------------------------------
# Sheet 6: "Confusion Matrix BHP GROUP (-0.6, 0.6, 0.6)"
# ----------------------------------------------------------------------
$wsCm5 = $wb.Worksheets.Item("Confusion Matrix BHP GROUP (-0.6, 0.6, 0.6)")
$wsCm5.Range("B3").Value = 4
$wsCm5.Range("C3").Value = 1551
$wsCm5.Range("D3").Value = 3
